# Add new "Partidos" (matches) rows for the 2025-11-01 session (Excel date serial 45962).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partidos")
$ws.Activate()

# r, fecha, jugador, equipo, posicion, goles, autogoles, arquero, goles_recibidos, tarjetas_amarillas, tarjetas_rojas, asistencias, Penales_Atajados
$rows = @(
    @(502, 45962, "Gember Marin Sarria", "Amarillo", "Arquero", 0, 0, $true, 4, 0, 0, 0, 0),
    @(503, 45962, "Jorge Gonzalez", "Azul", "Arquero", 0, 0, $true, 10, 0, 0, 0, 0),
    @(504, 45962, "Andres Tangarife", "Amarillo", "Delantero", 2, 0, $false, 0, 0, 0, 2, 0),
    @(505, 45962, "Cesar Augusto Estrada", "Amarillo", "Delantero", 2, 0, $false, 0, 0, 0, 0, 0),
    @(506, 45962, "Alexander Uribe", "Amarillo", "Mediocampista", 3, 0, $false, 0, 0, 0, 3, 0),
    @(507, 45962, "Armando Murillo", "Amarillo", "Defensa", 1, 0, $false, 0, 0, 0, 0, 0),
    @(508, 45962, "Julio Cesar Castaño", "Amarillo", "Mediocampista", 1, 0, $false, 0, 0, 0, 1, 0),
    @(509, 45962, "Andres Jurado", "Amarillo", "Delantero", 1, 0, $false, 0, 0, 0, 0, 0),
    @(510, 45962, "Andres Guerrero", "Amarillo", "Defensa", 0, 0, $false, 0, 0, 0, 2, 0),
    @(511, 45962, "Bryan Andres Burgos", "Amarillo", "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0),
    @(512, 45962, "Juan David Espinal", "Azul", "Mediocampista", 2, 0, $false, 0, 0, 0, 0, 0),
    @(513, 45962, "Arnul David Narvaez", "Azul", "Delantero", 1, 0, $false, 0, 0, 0, 1, 0),
    @(514, 45962, "David Fernando Velasco", "Azul", "Delantero", 1, 0, $false, 0, 0, 0, 0, 0),
    @(515, 45962, "Armando Vieras", "Azul", "Defensa", 0, 0, $false, 0, 0, 0, 2, 0),
    @(516, 45962, "Sebastian Giraldo", "Azul", "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]   # fecha
    $ws.Cells.Item($r, 2).Value = $row[2]   # jugador
    $ws.Cells.Item($r, 3).Value = $row[3]   # equipo
    $ws.Cells.Item($r, 4).Value = $row[4]   # posicion
    $ws.Cells.Item($r, 5).Value = $row[5]   # goles
    $ws.Cells.Item($r, 6).Value = $row[6]   # autogoles
    $ws.Cells.Item($r, 7).Value = $row[7]   # arquero
    $ws.Cells.Item($r, 8).Value = $row[8]   # goles_recibidos
    $ws.Cells.Item($r, 9).Value = $row[9]   # tarjetas_amarillas
    $ws.Cells.Item($r, 10).Value = $row[10] # tarjetas_rojas
    $ws.Cells.Item($r, 11).Value = $row[11] # asistencias
    $ws.Cells.Item($r, 12).Value = $row[12] # Penales_Atajados
}

# Re-apply the frozen header row and scroll/select so the newly added rows are in view,
# matching the author's last on-screen position after the edit.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$headerSel = $ws.Range("A2").Select()
$win.FreezePanes = $true
$win.ScrollRow = 500
$finalSel = $ws.Range("B517").Select()
